$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.768.62'
$ws.Range("E2").Value = '  -1.62%  '
$ws.Range("D3").Value = '2.675.88'
$ws.Range("E3").Value = '  -2.19%  '
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").Value = "'552.83"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.02%  '
$ws.Range("D6").Value = "'158.03"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.69%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("E8").Value = '  -0.64%  '
$ws.Range("E9").Value = '  -2.59%  '
$ws.Range("E10").Value = '  -2.58%  '
$ws.Range("E11").Value = '  -3.16%  '
$ws.Range("D12").Value = "'5.33"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -5.65%  '
$ws.Range("D13").Value = '3.149.26'
$ws.Range("E13").Value = '  -2.30%  '
$ws.Range("D14").Value = "'26.45"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.33%  '
$ws.Range("D15").Value = '62.670.48'
$ws.Range("E15").Value = '  -1.59%  '
$ws.Range("E16").Value = '  -1.35%  '
$ws.Range("D17").Value = '2.677.46'
$ws.Range("E17").Value = '  -2.39%  '
$ws.Range("E18").Value = '  -3.97%  '
$ws.Range("E19").Value = '  -3.06%  '
$ws.Range("D20").Value = "'344.59"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.42%  '
$ws.Range("E21").Value = '  -4.37%  '
$ws.Range("E22").Value = '  -0.09%  '
$ws.Range("D23").Value = "'0.507"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.75%  '
$ws.Range("D24").Value = "'63.18"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.51%  '
$ws.Range("E25").Value = '  -0.58%  '
$ws.Range("E26").Value = '  -0.25%  '
$ws.Range("D27").Value = "'8.21"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.55%  '
$ws.Range("E28").Value = '  +9.53%  '
$ws.Range("D29").Value = '0.0₃0853'
$ws.Range("E29").Value = '  -5.30%  '
$ws.Range("E30").Value = '  +0.84%  '
$ws.Range("E31").Value = '  -1.03%  '
$ws.Range("D32").Value = "'163.06"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.37%  '
$ws.Range("D33").Value = "'4.91"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.06%  '
$ws.Range("E34").Value = '  +0.74%  '
$ws.Range("E35").Value = '  +0.01%  '
$ws.Range("D36").Value = "'19.46"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.91%  '
$ws.Range("E37").Value = '  -0.76%  '
$ws.Range("D38").Value = "'351.98"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.72%  '
$ws.Range("D39").Value = "'0.949"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -3.23%  '
$ws.Range("D40").Value = "'6.22"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.76%  '
$ws.Range("E41").Value = '  -1.95%  '
$ws.Range("D42").Value = "'38.36"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.10%  '
$ws.Range("D43").Value = "'20.92"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.63%  '
$ws.Range("D44").Value = "'20.21"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -3.72%  '
$ws.Range("D45").Value = "'0.615"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.09%  '
$ws.Range("E46").Value = '  -3.37%  '
$ws.Range("D47").Value = "'0.998"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.17%  '
$ws.Range("E48").Value = '  -0.41%  '
$ws.Range("B49").Value = 'VeChain'
$ws.Range("C49").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D49").Value = "'0.0243"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.61%  '
$ws.Range("B50").Value = 'Stellar'
$ws.Range("C50").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D50").Value = "'0.0972"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.75%  '
$ws.Range("D51").Value = "'128.74"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -4.59%  '

Write-Output "Done applying crypto updates"
